# Auto-generated Excel COM-interop script applying the Durandal_Profits.xlsx
# profit-recalculation update across all eight crafting-job sheets.
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 3000
$ws.Range("I69").Value = 1750
$ws.Range("J69").Value = 3625
$ws.Range("K69").Value = 5250
$ws.Range("L69").Value = 10875
$ws.Range("M69").Value = -4376
$ws.Range("N69").Value = -12623
$ws.Range("H72").Value = 3000
$ws.Range("I72").Value = 1750
$ws.Range("J72").Value = 3625
$ws.Range("K72").Value = 15750
$ws.Range("L72").Value = 32625
$ws.Range("M72").Value = -11382
$ws.Range("N72").Value = -41361
$ws.Range("H74").Value = 3900
$ws.Range("I74").Value = 3828.5715
$ws.Range("J74").Value = 4400
$ws.Range("K74").Value = 3828.5715
$ws.Range("L74").Value = 4400
$ws.Range("M74").Value = -2892.5715
$ws.Range("N74").Value = -6272
$ws.Range("H77").Value = 3900
$ws.Range("I77").Value = 3828.5715
$ws.Range("J77").Value = 4400
$ws.Range("K77").Value = 19142.8575
$ws.Range("L77").Value = 22000
$ws.Range("M77").Value = -14462.8575
$ws.Range("N77").Value = -31360
$ws.Range("H138").Value = 3800.6843
$ws.Range("I138").Value = 1782.8
$ws.Range("J138").Value = 5377.1562
$ws.Range("K138").Value = 5348.4
$ws.Range("L138").Value = 16131.4686
$ws.Range("M138").Value = -208.3999999999996
$ws.Range("N138").Value = -26411.4686

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 350832.22
$ws.Range("I32").Value = 2017.7838
$ws.Range("J32").Value = 3577365.8
$ws.Range("K32").Value = 2017.7838
$ws.Range("L32").Value = 3577365.8
$ws.Range("M32").Value = -1730.7838
$ws.Range("N32").Value = -3577939.8
$ws.Range("H122").Value = 45687.957
$ws.Range("I122").Value = 2465.45
$ws.Range("J122").Value = 333838
$ws.Range("K122").Value = 7396.349999999999
$ws.Range("L122").Value = 1001514
$ws.Range("M122").Value = -4946.349999999999
$ws.Range("N122").Value = -1006414
$ws.Range("H132").Value = 24881.072
$ws.Range("I132").Value = 650.8788
$ws.Range("J132").Value = 113725.11
$ws.Range("K132").Value = 1952.6364
$ws.Range("L132").Value = 341175.33
$ws.Range("M132").Value = 577.3636000000001
$ws.Range("N132").Value = -346235.33

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 31252582
$ws.Range("I86").Value = 43480644
$ws.Range("J86").Value = 3090.4443
$ws.Range("K86").Value = 43480644
$ws.Range("L86").Value = 3090.4443
$ws.Range("M86").Value = -43479521
$ws.Range("N86").Value = -5336.4443
$ws.Range("H89").Value = 31252582
$ws.Range("I89").Value = 43480644
$ws.Range("J89").Value = 3090.4443
$ws.Range("K89").Value = 217403220
$ws.Range("L89").Value = 15452.2215
$ws.Range("M89").Value = -217397604
$ws.Range("N89").Value = -26684.2215
$ws.Range("H94").Value = 824.8
$ws.Range("I94").Value = 735.8461
$ws.Range("J94").Value = 990
$ws.Range("K94").Value = 735.8461
$ws.Range("L94").Value = 990
$ws.Range("M94").Value = -284.8461
$ws.Range("N94").Value = -1892
$ws.Range("H134").Value = 8127.5
$ws.Range("I134").Value = 2437.5833
$ws.Range("J134").Value = 76406.5
$ws.Range("K134").Value = 7312.749899999999
$ws.Range("L134").Value = 229219.5
$ws.Range("M134").Value = -4777.749899999999
$ws.Range("N134").Value = -234289.5

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4718805
$ws.Range("I31").Value = 6411661.5
$ws.Range("J31").Value = 2989.2144
$ws.Range("K31").Value = 6411661.5
$ws.Range("L31").Value = 2989.2144
$ws.Range("M31").Value = -6411366.5
$ws.Range("N31").Value = -3579.2144
$ws.Range("H34").Value = 4718805
$ws.Range("I34").Value = 6411661.5
$ws.Range("J34").Value = 2989.2144
$ws.Range("K34").Value = 6411661.5
$ws.Range("L34").Value = 2989.2144
$ws.Range("M34").Value = -6411459.5
$ws.Range("N34").Value = -3393.2144
$ws.Range("H58").Value = 1126.375
$ws.Range("I58").Value = 1238.5454
$ws.Range("J58").Value = 879.6
$ws.Range("K58").Value = 1238.5454
$ws.Range("L58").Value = 879.6
$ws.Range("M58").Value = -1035.5454
$ws.Range("N58").Value = -1285.6
$ws.Range("H86").Value = 55631.35
$ws.Range("I86").Value = 21500
$ws.Range("J86").Value = 69852.75
$ws.Range("K86").Value = 21500
$ws.Range("L86").Value = 69852.75
$ws.Range("M86").Value = -20377
$ws.Range("N86").Value = -72098.75
$ws.Range("H89").Value = 55631.35
$ws.Range("I89").Value = 21500
$ws.Range("J89").Value = 69852.75
$ws.Range("K89").Value = 107500
$ws.Range("L89").Value = 349263.75
$ws.Range("M89").Value = -101884
$ws.Range("N89").Value = -360495.75
$ws.Range("H136").Value = 1126.375
$ws.Range("I136").Value = 1238.5454
$ws.Range("J136").Value = 879.6
$ws.Range("K136").Value = 3715.6362
$ws.Range("L136").Value = 2638.8
$ws.Range("M136").Value = -1165.6362
$ws.Range("N136").Value = -7738.8

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 642.46155
$ws.Range("I109").Value = 316.8889
$ws.Range("J109").Value = 1375
$ws.Range("K109").Value = 950.6667
$ws.Range("L109").Value = 4125
$ws.Range("M109").Value = 89.33330000000001
$ws.Range("N109").Value = -6205
$ws.Range("H120").Value = 8994.286
$ws.Range("I120").Value = 4500
$ws.Range("J120").Value = 14986.667
$ws.Range("K120").Value = 13500
$ws.Range("L120").Value = 44960.001
$ws.Range("M120").Value = -8662
$ws.Range("N120").Value = -54636.001
$ws.Range("H128").Value = 0
$ws.Range("I128").Value = 0
$ws.Range("K128").Value = 0
$ws.Range("M128").ClearContents()
$ws.Range("H131").Value = 26317160
$ws.Range("J131").Value = 55557410
$ws.Range("L131").Value = 166672230
$ws.Range("N131").Value = -166682310

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H102").Value = 1894.4375
$ws.Range("I102").Value = 1485.5385
$ws.Range("J102").Value = 3666.3333
$ws.Range("K102").Value = 1485.5385
$ws.Range("L102").Value = 3666.3333
$ws.Range("M102").Value = 136.4614999999999
$ws.Range("N102").Value = -6910.3333

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2083.3333
$ws.Range("I82").Value = 1750
$ws.Range("J82").Value = 2250
$ws.Range("K82").Value = 1750
$ws.Range("L82").Value = 2250
$ws.Range("M82").Value = -1389
$ws.Range("N82").Value = -2972
$ws.Range("H85").Value = 2083.3333
$ws.Range("I85").Value = 1750
$ws.Range("J85").Value = 2250
$ws.Range("K85").Value = 1750
$ws.Range("L85").Value = 2250
$ws.Range("M85").Value = -502
$ws.Range("N85").Value = -4746
$ws.Range("H93").Value = 1048.8182
$ws.Range("I93").Value = 737.56525
$ws.Range("J93").Value = 1764.7
$ws.Range("K93").Value = 737.56525
$ws.Range("L93").Value = 1764.7
$ws.Range("M93").Value = 510.43475
$ws.Range("N93").Value = -4260.7

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3413.1
$ws.Range("I81").Value = 991
$ws.Range("J81").Value = 3682.2222
$ws.Range("K81").Value = 1982
$ws.Range("L81").Value = 7364.4444
$ws.Range("M81").Value = -921
$ws.Range("N81").Value = -9486.4444
$ws.Range("H84").Value = 3413.1
$ws.Range("I84").Value = 991
$ws.Range("J84").Value = 3682.2222
$ws.Range("K84").Value = 9910
$ws.Range("L84").Value = 36822.222
$ws.Range("M84").Value = -4606
$ws.Range("N84").Value = -47430.222
$ws.Range("H135").Value = 67688.8
$ws.Range("J135").Value = 67688.8
$ws.Range("L135").Value = 67688.8
$ws.Range("N135").Value = -77828.8
$ws.Range("H136").Value = 32104.719
$ws.Range("I136").Value = 42310.043
$ws.Range("J136").Value = 1488.75
$ws.Range("K136").Value = 126930.129
$ws.Range("L136").Value = 4466.25
$ws.Range("M136").Value = -124380.129
$ws.Range("N136").Value = -9566.25

